# Update column G (K) values per row, per commit:
# "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gUpdates = @(
    @{Row=2; Value=1},
    @{Row=3; Value=0},
    @{Row=4; Value=0},
    @{Row=5; Value=0},
    @{Row=6; Value=1},
    @{Row=7; Value=1},
    @{Row=8; Value=1},
    @{Row=9; Value=0},
    @{Row=10; Value=2},
    @{Row=11; Value=0},
    @{Row=12; Value=1},
    @{Row=14; Value=0},
    @{Row=15; Value=1},
    @{Row=16; Value=3},
    @{Row=17; Value=0},
    @{Row=18; Value=0},
    @{Row=19; Value=0},
    @{Row=20; Value=1},
    @{Row=21; Value=3},
    @{Row=22; Value=0},
    @{Row=23; Value=1},
    @{Row=24; Value=0},
    @{Row=25; Value=2},
    @{Row=26; Value=0},
    @{Row=27; Value=2},
    @{Row=28; Value=1},
    @{Row=29; Value=2},
    @{Row=30; Value=0},
    @{Row=31; Value=0},
    @{Row=32; Value=2},
    @{Row=33; Value=0},
    @{Row=34; Value=0},
    @{Row=35; Value=0},
    @{Row=36; Value=0},
    @{Row=37; Value=2},
    @{Row=38; Value=3},
    @{Row=39; Value=1},
    @{Row=40; Value=1},
    @{Row=41; Value=0},
    @{Row=42; Value=0},
    @{Row=43; Value=1},
    @{Row=44; Value=1},
    @{Row=45; Value=0},
    @{Row=46; Value=0},
    @{Row=47; Value=0},
    @{Row=48; Value=1},
    @{Row=49; Value=1},
    @{Row=50; Value=2},
    @{Row=51; Value=0},
    @{Row=52; Value=1},
    @{Row=53; Value=1},
    @{Row=54; Value=0},
    @{Row=55; Value=2},
    @{Row=56; Value=0},
    @{Row=57; Value=1},
    @{Row=58; Value=0},
    @{Row=59; Value=0},
    @{Row=60; Value=1},
    @{Row=61; Value=1},
    @{Row=62; Value=0},
    @{Row=63; Value=0},
    @{Row=64; Value=3},
    @{Row=65; Value=0},
    @{Row=66; Value=0},
    @{Row=67; Value=0},
    @{Row=68; Value=0},
    @{Row=69; Value=1},
    @{Row=70; Value=0},
    @{Row=71; Value=0},
    @{Row=72; Value=1},
    @{Row=73; Value=0},
    @{Row=74; Value=0},
    @{Row=75; Value=0},
    @{Row=76; Value=0},
    @{Row=77; Value=1},
    @{Row=78; Value=0},
    @{Row=79; Value=3},
    @{Row=80; Value=0},
    @{Row=81; Value=1},
    @{Row=82; Value=1},
    @{Row=83; Value=2},
    @{Row=84; Value=1},
    @{Row=85; Value=0},
    @{Row=87; Value=1},
    @{Row=89; Value=1},
    @{Row=91; Value=1},
    @{Row=92; Value=1}
)

foreach ($u in $gUpdates) {
    $ws.Cells.Item($u.Row, 7).Value = $u.Value
}
